# kfold para ML - DONE
# Fill in the k-fold cross-validation performance measurements on the
# "resultados" sheet (D4:D7), give them a 6-decimal numeric format, make
# "resultados" the active/selected sheet+range, and drop the underline
# formatting that used to sit (unused) on D5.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)          # resultados

# --- resultados: new performance numbers -----------------------------
$ws1.Activate()

# D5 already carried a (unused/underlined) font style from before; strip
# that off so the new values render with the normal font.
$ws1.Range("D5").Font.Underline = $false

$ws1.Range("D4").Value = 0.056794900000000002
$ws1.Range("D5").Value = 0.0055655879999999998
$ws1.Range("D6").Value = 0.074602870000000002
$ws1.Range("D7").Value = 0.3723207

$ws1.Range("D4:D7").NumberFormat = "0.000000"

# Match the saved selection state: resultados active, D4:D7 selected.
$ws1.Range("D4:D7").Select()
